$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45436

$ws.Range("D14").Value = 41.9
$ws.Range("D15").Value = 47.8
$ws.Range("D16").Value = 57.85
$ws.Range("D17").Value = 73.7
$ws.Range("D18").Value = 147.2
$ws.Range("D19").Value = 221
$ws.Range("D20").Value = 344
$ws.Range("D21").Value = 475
